$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (A2:D12, header is row 1) need to be re-sorted in ascending
# order by column A ("time (s)"). Read all rows into memory, sort them by
# the first column, then write the sorted rows back in place.

$firstRow = 2
$lastRow = 12
$firstCol = 1
$lastCol = 4

$rows = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value()
    }
    $rows += ,$rowVals
}

$sortedRows = $rows | Sort-Object -Property { $_[0] }

for ($i = 0; $i -lt $sortedRows.Count; $i++) {
    $r = $firstRow + $i
    $rowVals = $sortedRows[$i]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - $firstCol]
    }
}
